$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header F1 with same style as header row (copy format from E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = 'Chinese'

# Build a helper cell with the Menlo font used for F2:F40, then copy that format over
$helper = $ws.Range("H1")
$helper.Value = "x"
$helper.Font.Name = "Menlo"
$helper.Font.Size = 12
$helper.Font.Color = 0
$helper.Copy() | Out-Null
$ws.Range("F2:F40").PasteSpecial(-4122) | Out-Null
$helper.Clear() | Out-Null

# Fill in the Chinese translations
$ws.Range("F2").Value = '必填'
$ws.Range("F3").Value = '個人所得稅計算器'
$ws.Range("F4").Value = '收入'
$ws.Range("F5").Value = '輸入收入'
$ws.Range("F6").Value = '受扶養人'
$ws.Range("F7").Value = '輸入受扶養人'
$ws.Range("F8").Value = '繳費等級'
$ws.Range("F9").Value = '以正式工資計算'
$ws.Range("F10").Value = '其他'
$ws.Range("F11").Value = '繳費金額'
$ws.Range("F12").Value = '至少 {0}₫'
$ws.Range("F13").Value = '目標工資類型'
$ws.Range("F14").Value = '淨收入'
$ws.Range("F15").Value = '總收入'
$ws.Range("F16").Value = '計算至 {0}'
$ws.Range("F17").Value = '薪資和稅務詳情'
$ws.Range("F18").Value = '總收入為 {0}₫'
$ws.Range("F19").Value = '已繳保險費為 {0}₫'
$ws.Range("F20").Value = '稅前收入為 {0}₫'
$ws.Range("F21").Value = '應稅所得為 {0}₫'
$ws.Range("F22").Value = '稅額為 {0}₫'
$ws.Range("F23").Value = '淨收入為 {0}₫'
$ws.Range("F24").Value = '設定'
$ws.Range("F25").Value = '語言'
$ws.Range("F26").Value = '深色模式'
$ws.Range("F27").Value = '關閉'
$ws.Range("F28").Value = '個人扣除額'
$ws.Range("F29").Value = '受扶養人扣除額'
$ws.Range("F30").Value = '保險費率'
$ws.Range("F31").Value = '最低保險基數'
$ws.Range("F32").Value = '個人所得稅政策'
$ws.Range("F33").Value = '04/2007/QH12 號 法律'
$ws.Range("F34").Value = '第 954/2020/UBTVQH14 號決議'
$ws.Range("F35").Value = '第 110/2025/UBTVQH15 號決議'
$ws.Range("F36").Value = '稅收政策詳情'
$ws.Range("F37").Value = '稅級'
$ws.Range("F38").Value = '稅階'
$ws.Range("F39").Value = '稅率'
$ws.Range("F40").Value = '關閉'

# Set column F width to match target (36.5 excel units)
$ws.Columns("F").ColumnWidth = 35 + 2/3

# Update sheet view / selection to match target state
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("A2").Select() | Out-Null
$ws.Range("F6").Select() | Out-Null
